# Klepelere kadar(klepeler dahil) kurulum ayarları kısmı tamamlandı.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("EkranDatabase")
$ws2 = $wb.Worksheets.Item("ServerDBKurulum")

# --- EkranDatabase: add two new rows (17, 18) ---
$ws1.Range("A17").Value = 15
$ws1.Range("B17").Value = 16
$ws1.Range("C17").Value = "Klepe Harita Onay"
$ws1.Range("D17").Value = "klepeHaritalar"
$ws1.Range("E17").Value = 0
$ws1.Range("F17").Value = 0

$ws1.Range("A18").Value = 16
$ws1.Range("B18").Value = 17
$ws1.Range("C18").Value = "Klepe ve Çıkış No Onay"
$ws1.Range("D18").Value = "klepeNo'lar"
$ws1.Range("E18").Value = "cikisNo'lar"
$ws1.Range("F18").Value = 0

# --- ServerDBKurulum: add one new row (22) ---
$ws2.Range("A22").Value = 20
$ws2.Range("B22").Value = 21
$ws2.Range("C22").Value = "Klepe No ve Çıkış No Durum(ok veya null)"
$ws2.Range("D22").Value = "fanNo'lar"
$ws2.Range("E22").Value = "çıkışNo'lar"
$ws2.Range("F22").Value = 0

# Column C on ServerDBKurulum grew to fit the new, longer text (no longer best-fit flagged)
$ws2.Columns.Item(3).ColumnWidth = 37.6

# Update selections on both sheets to match the saved view state
$ws1.Range("E25").Select()

# ServerDBKurulum becomes the active/visible tab
$ws2.Activate()
$ws2.Range("D22").Select()
